$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 75
$ws.Range("I4").Value = 75
$ws.Range("K4").Value = 75
$ws.Range("M4").Value = 39
$ws.Range("H15").Value = 1817.174
$ws.Range("I15").Value = 1817.174
$ws.Range("K15").Value = 5451.522
$ws.Range("M15").Value = -5282.522
$ws.Range("H33").Value = 175.83333
$ws.Range("I33").Value = 187.8125
$ws.Range("K33").Value = 187.8125
$ws.Range("M33").Value = 41.1875
$ws.Range("H61").Value = 909
$ws.Range("I61").Value = 909
$ws.Range("K61").Value = 2727
$ws.Range("M61").Value = -2555
$ws.Range("H115").Value = 1019.8
$ws.Range("I115").Value = 1019.8
$ws.Range("K115").Value = 3059.4
$ws.Range("M115").Value = -1492.4
$ws.Range("H137").Value = 1138.8334
$ws.Range("I137").Value = 940.375
$ws.Range("K137").Value = 2821.125
$ws.Range("M137").Value = -271.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1361.6875
$ws.Range("I2").Value = 1389.1333
$ws.Range("K2").Value = 1389.1333
$ws.Range("M2").Value = -1276.1333
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H45").Value = 1316.25
$ws.Range("I45").Value = 1850
$ws.Range("J45").Value = 782.5
$ws.Range("K45").Value = 1850
$ws.Range("L45").Value = 782.5
$ws.Range("M45").Value = -1473
$ws.Range("N45").Value = -1536.5
$ws.Range("H61").Value = 2129.2222
$ws.Range("I61").Value = 2129.2222
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2129.2222
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1917.2222
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 16332.125
$ws.Range("I74").Value = 14378.143
$ws.Range("K74").Value = 14378.143
$ws.Range("M74").Value = -13504.143
$ws.Range("H77").Value = 16332.125
$ws.Range("I77").Value = 14378.143
$ws.Range("K77").Value = 71890.715
$ws.Range("M77").Value = -67522.715
$ws.Range("H116").Value = 1361.6875
$ws.Range("I116").Value = 1389.1333
$ws.Range("K116").Value = 1389.1333
$ws.Range("M116").Value = 904.8667
$ws.Range("H132").Value = 1179.025
$ws.Range("I132").Value = 907.56525
$ws.Range("K132").Value = 2722.69575
$ws.Range("M132").Value = -192.6957499999999
$ws.Range("H136").Value = 2129.2222
$ws.Range("I136").Value = 2129.2222
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6387.6666
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3837.6666
$ws.Range("N136").ClearContents()
$ws.Range("H139").Value = 85000
$ws.Range("J139").Value = 85000
$ws.Range("L139").Value = 85000
$ws.Range("N139").Value = -95280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1361.6875
$ws.Range("I3").Value = 1389.1333
$ws.Range("K3").Value = 1389.1333
$ws.Range("M3").Value = -1275.1333
$ws.Range("H80").Value = 1498.5
$ws.Range("J80").Value = 1501.1428
$ws.Range("L80").Value = 1501.1428
$ws.Range("N80").Value = -3497.1428
$ws.Range("H82").Value = 44030.832
$ws.Range("I82").Value = 4200
$ws.Range("J82").Value = 51997
$ws.Range("K82").Value = 4200
$ws.Range("L82").Value = 51997
$ws.Range("M82").Value = -3817
$ws.Range("N82").Value = -52763
$ws.Range("H83").Value = 1498.5
$ws.Range("J83").Value = 1501.1428
$ws.Range("L83").Value = 7505.714
$ws.Range("N83").Value = -17489.714
$ws.Range("H85").Value = 44030.832
$ws.Range("I85").Value = 4200
$ws.Range("J85").Value = 51997
$ws.Range("K85").Value = 4200
$ws.Range("L85").Value = 51997
$ws.Range("M85").Value = -2874
$ws.Range("N85").Value = -54649
$ws.Range("H86").Value = 2282.9167
$ws.Range("I86").Value = 2043.8889
$ws.Range("K86").Value = 2043.8889
$ws.Range("M86").Value = -920.8888999999999
$ws.Range("H89").Value = 2282.9167
$ws.Range("I89").Value = 2043.8889
$ws.Range("K89").Value = 10219.4445
$ws.Range("M89").Value = -4603.4445
$ws.Range("H94").Value = 709
$ws.Range("I94").Value = 685.75
$ws.Range("K94").Value = 685.75
$ws.Range("M94").Value = -234.75
$ws.Range("H134").Value = 2572.9375
$ws.Range("I134").Value = 2411.1333
$ws.Range("K134").Value = 7233.3999
$ws.Range("M134").Value = -4698.3999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3874.6924
$ws.Range("I132").Value = 2964.5
$ws.Range("K132").Value = 8893.5
$ws.Range("M132").Value = -6363.5
$ws.Range("H141").Value = 51788
$ws.Range("J141").Value = 59985
$ws.Range("L141").Value = 59985
$ws.Range("N141").Value = -70345

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 5500
$ws.Range("I70").Value = 1000
$ws.Range("J70").Value = 10000
$ws.Range("K70").Value = 3000
$ws.Range("L70").Value = 30000
$ws.Range("M70").Value = -2685
$ws.Range("N70").Value = -30630
$ws.Range("H73").Value = 5500
$ws.Range("I73").Value = 1000
$ws.Range("J73").Value = 10000
$ws.Range("K73").Value = 3000
$ws.Range("L73").Value = 30000
$ws.Range("M73").Value = -1908
$ws.Range("N73").Value = -32184
$ws.Range("H121").Value = 1939.4445
$ws.Range("I121").Value = 998
$ws.Range("J121").Value = 2208.4285
$ws.Range("K121").Value = 2994
$ws.Range("L121").Value = 6625.2855
$ws.Range("M121").Value = -1684
$ws.Range("N121").Value = -9245.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1075.6666
$ws.Range("I126").Value = 1056.5
$ws.Range("J126").Value = 1114
$ws.Range("K126").Value = 3169.5
$ws.Range("L126").Value = 3342
$ws.Range("M126").Value = -699.5
$ws.Range("N126").Value = -8282
$ws.Range("H132").Value = 3436.077
$ws.Range("I132").Value = 3192.6667
$ws.Range("J132").Value = 3509.1
$ws.Range("K132").Value = 9578.000100000001
$ws.Range("L132").Value = 10527.3
$ws.Range("M132").Value = -7048.000100000001
$ws.Range("N132").Value = -15587.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2283.625
$ws.Range("I7").Value = 2283.625
$ws.Range("K7").Value = 2283.625
$ws.Range("M7").Value = -2171.625
$ws.Range("H82").Value = 2999.5
$ws.Range("I82").Value = 2999.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2999.5
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -2638.5
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 2999.5
$ws.Range("I85").Value = 2999.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2999.5
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -1751.5
$ws.Range("N85").ClearContents()
$ws.Range("H126").Value = 2283.625
$ws.Range("I126").Value = 2283.625
$ws.Range("K126").Value = 6850.875
$ws.Range("M126").Value = -4380.875
$ws.Range("H132").Value = 5571.6
$ws.Range("I132").Value = 4398.5
$ws.Range("J132").Value = 5998.1816
$ws.Range("K132").Value = 13195.5
$ws.Range("L132").Value = 17994.5448
$ws.Range("M132").Value = -10665.5
$ws.Range("N132").Value = -23054.5448
